# genSymbolValsWithPos regenerated its symbol/position rows after the
# isNeedGen flag was threaded through - row order (A2:A25) changed while
# the underlying symbol/reel value set stayed the same. Rewrite the data
# rows (headers in row 1 and totals in row 26 are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,6
$arr[0,0]=902;$arr[0,1]=1;$arr[0,2]=0;$arr[0,3]=0;$arr[0,4]=0;$arr[0,5]=0
$arr[1,0]=1001;$arr[1,1]=18;$arr[1,2]=30;$arr[1,3]=75;$arr[1,4]=60;$arr[1,5]=72
$arr[2,0]=701;$arr[2,1]=3;$arr[2,2]=90;$arr[2,3]=45;$arr[2,4]=97;$arr[2,5]=15
$arr[3,0]=201;$arr[3,1]=9;$arr[3,2]=30;$arr[3,3]=15;$arr[3,4]=45;$arr[3,5]=30
$arr[4,0]=501;$arr[4,1]=9;$arr[4,2]=52;$arr[4,3]=30;$arr[4,4]=75;$arr[4,5]=45
$arr[5,0]=401;$arr[5,1]=9;$arr[5,2]=48;$arr[5,3]=67;$arr[5,4]=75;$arr[5,5]=45
$arr[6,0]=601;$arr[6,1]=9;$arr[6,2]=60;$arr[6,3]=67;$arr[6,4]=60;$arr[6,5]=42
$arr[7,0]=101;$arr[7,1]=9;$arr[7,2]=30;$arr[7,3]=15;$arr[7,4]=60;$arr[7,5]=15
$arr[8,0]=801;$arr[8,1]=3;$arr[8,2]=67;$arr[8,3]=65;$arr[8,4]=52;$arr[8,5]=45
$arr[9,0]=1203;$arr[9,1]=3;$arr[9,2]=15;$arr[9,3]=15;$arr[9,4]=15;$arr[9,5]=15
$arr[10,0]=901;$arr[10,1]=16;$arr[10,2]=15;$arr[10,3]=45;$arr[10,4]=60;$arr[10,5]=60
$arr[11,0]=301;$arr[11,1]=6;$arr[11,2]=45;$arr[11,3]=30;$arr[11,4]=60;$arr[11,5]=45
$arr[12,0]=1201;$arr[12,1]=2;$arr[12,2]=10;$arr[12,3]=10;$arr[12,4]=10;$arr[12,5]=10
$arr[13,0]=1202;$arr[13,1]=2;$arr[13,2]=10;$arr[13,3]=10;$arr[13,4]=10;$arr[13,5]=10
$arr[14,0]=2;$arr[14,1]=0;$arr[14,2]=2;$arr[14,3]=2;$arr[14,4]=2;$arr[14,5]=2
$arr[15,0]=802;$arr[15,1]=0;$arr[15,2]=4;$arr[15,3]=5;$arr[15,4]=4;$arr[15,5]=0
$arr[16,0]=502;$arr[16,1]=0;$arr[16,2]=4;$arr[16,3]=0;$arr[16,4]=0;$arr[16,5]=0
$arr[17,0]=1101;$arr[17,1]=0;$arr[17,2]=15;$arr[17,3]=30;$arr[17,4]=30;$arr[17,5]=0
$arr[18,0]=1;$arr[18,1]=0;$arr[18,2]=2;$arr[18,3]=2;$arr[18,4]=2;$arr[18,5]=2
$arr[19,0]=3;$arr[19,1]=0;$arr[19,2]=3;$arr[19,3]=3;$arr[19,4]=3;$arr[19,5]=3
$arr[20,0]=402;$arr[20,1]=0;$arr[20,2]=0;$arr[20,3]=4;$arr[20,4]=0;$arr[20,5]=0
$arr[21,0]=602;$arr[21,1]=0;$arr[21,2]=0;$arr[21,3]=4;$arr[21,4]=0;$arr[21,5]=9
$arr[22,0]=702;$arr[22,1]=0;$arr[22,2]=0;$arr[22,3]=0;$arr[22,4]=4;$arr[22,5]=0
$arr[23,0]=1002;$arr[23,1]=0;$arr[23,2]=0;$arr[23,3]=0;$arr[23,4]=0;$arr[23,5]=9
$ws.Range("A2:F25").Value = $arr
